# Natmi following Dr Hou advice
#
# The NATMI Clcf1 -> Cntfr edge table is regenerated with an additional
# "ECs" sending cluster (on top of the existing FAPs / M2 / sCs clusters),
# each evaluated against the same two target clusters (FAPs, sCs). This
# grows the data block from 6 rows (3 senders x 2 targets) to 8 rows
# (4 senders x 2 targets) and refreshes every derived statistic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Clcf1"
$ws.Range("C2").Value = "Cntfr"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.679012
$ws.Range("H2").Value = 5.037036000000001
$ws.Range("I2").Value = 0.1178149724053671
$ws.Range("J2").Value = 0.1178149724053671
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.666771333333333
$ws.Range("N2").Value = 17.000314
$ws.Range("O2").Value = 0.9648745059153377
$ws.Range("P2").Value = 0.9648745059153376
$ws.Range("Q2").Value = 9.514577069922668
$ws.Range("R2").Value = 85.63119362930401
$ws.Range("S2").Value = 0.1136766632890578
$ws.Range("T2").Value = 0.1136766632890578

# Row 3: ECs -> sCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Clcf1"
$ws.Range("C3").Value = "Cntfr"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.679012
$ws.Range("H3").Value = 5.037036000000001
$ws.Range("I3").Value = 0.1178149724053671
$ws.Range("J3").Value = 0.1178149724053671
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2062943333333333
$ws.Range("N3").Value = 0.618883
$ws.Range("O3").Value = 0.03512549408466232
$ws.Range("P3").Value = 0.03512549408466231
$ws.Range("Q3").Value = 0.3463706611986667
$ws.Range("R3").Value = 3.117335950788
$ws.Range("S3").Value = 0.004138309116309378
$ws.Range("T3").Value = 0.004138309116309377

# Row 4: FAPs -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Clcf1"
$ws.Range("C4").Value = "Cntfr"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.197979
$ws.Range("H4").Value = 9.593937
$ws.Range("I4").Value = 0.2243997110431275
$ws.Range("J4").Value = 0.2243997110431275
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.666771333333333
$ws.Range("N4").Value = 17.000314
$ws.Range("O4").Value = 0.9648745059153377
$ws.Range("P4").Value = 0.9648745059153376
$ws.Range("Q4").Value = 18.122215721802
$ws.Range("R4").Value = 163.099941496218
$ws.Range("S4").Value = 0.2165175603202822
$ws.Range("T4").Value = 0.2165175603202822

# Row 5: FAPs -> sCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Clcf1"
$ws.Range("C5").Value = "Cntfr"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.197979
$ws.Range("H5").Value = 9.593937
$ws.Range("I5").Value = 0.2243997110431275
$ws.Range("J5").Value = 0.2243997110431275
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2062943333333333
$ws.Range("N5").Value = 0.618883
$ws.Range("O5").Value = 0.03512549408466232
$ws.Range("P5").Value = 0.03512549408466231
$ws.Range("Q5").Value = 0.659724945819
$ws.Range("R5").Value = 5.937524512371
$ws.Range("S5").Value = 0.007882150722845309
$ws.Range("T5").Value = 0.007882150722845307

# Row 6: M2 -> FAPs
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Clcf1"
$ws.Range("C6").Value = "Cntfr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.919382666666667
$ws.Range("H6").Value = 5.758148
$ws.Range("I6").Value = 0.134681596027112
$ws.Range("J6").Value = 0.134681596027112
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.666771333333333
$ws.Range("N6").Value = 17.000314
$ws.Range("O6").Value = 0.9648745059153377
$ws.Range("P6").Value = 0.9648745059153376
$ws.Range("Q6").Value = 10.87670267316355
$ws.Range("R6").Value = 97.890324058472
$ws.Range("S6").Value = 0.1299508384225488
$ws.Range("T6").Value = 0.1299508384225488

# Row 7: M2 -> sCs
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Clcf1"
$ws.Range("C7").Value = "Cntfr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.919382666666667
$ws.Range("H7").Value = 5.758148
$ws.Range("I7").Value = 0.134681596027112
$ws.Range("J7").Value = 0.134681596027112
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2062943333333333
$ws.Range("N7").Value = 0.618883
$ws.Range("O7").Value = 0.03512549408466232
$ws.Range("P7").Value = 0.03512549408466231
$ws.Range("Q7").Value = 0.3959577676315555
$ws.Range("R7").Value = 3.563619908684
$ws.Range("S7").Value = 0.004730757604563201
$ws.Range("T7").Value = 0.0047307576045632

# Row 8: sCs -> FAPs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Clcf1"
$ws.Range("C8").Value = "Cntfr"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.454888
$ws.Range("H8").Value = 22.364664
$ws.Range("I8").Value = 0.5231037205243934
$ws.Range("J8").Value = 0.5231037205243932
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.666771333333333
$ws.Range("N8").Value = 17.000314
$ws.Range("O8").Value = 0.9648745059153377
$ws.Range("P8").Value = 0.9648745059153376
$ws.Range("Q8").Value = 42.24514561161067
$ws.Range("R8").Value = 380.206310504496
$ws.Range("S8").Value = 0.504729443883449
$ws.Range("T8").Value = 0.5047294438834488

# Row 9: sCs -> sCs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Clcf1"
$ws.Range("C9").Value = "Cntfr"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.454888
$ws.Range("H9").Value = 22.364664
$ws.Range("I9").Value = 0.5231037205243934
$ws.Range("J9").Value = 0.5231037205243932
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2062943333333333
$ws.Range("N9").Value = 0.618883
$ws.Range("O9").Value = 0.03512549408466232
$ws.Range("P9").Value = 0.03512549408466231
$ws.Range("Q9").Value = 1.537901150034667
$ws.Range("R9").Value = 13.841110350312
$ws.Range("S9").Value = 0.01837427664094443
$ws.Range("T9").Value = 0.01837427664094442
